$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# The plan entry "20190103 / 小区总水、电表费用结算并过户" (row 4 on Sheet1)
# has been completed / superseded - delete the whole row, shifting the
# remaining plan rows up by one.
$ws1.Rows(4).Delete()

# Restore the on-sheet selections to match where the user ended up working.
$ws1.Range("B19").Select()
$ws2.Range("C2:C24").Select()
